# Revert "StfUtils Unit Tests TestCaseStepFilePathUtils: Added Tc4331 and Tc43311"
#
# The previous commit had added a new "Tc4331"/"Tc43311" test-case block
# (rows 41-43) to the TestCaseStepFilePathUtils sheet and a brand new
# worksheet ("More tests ") containing some scratch notes. This reverts
# both: the Tc4331/Tc43311 rows are removed, and the "More tests " sheet
# content is folded back into the bottom of the TestCaseStepFilePathUtils
# sheet (rows 42-46), after which the now-empty "More tests " sheet is
# deleted.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$excel.DisplayAlerts = $false

# Remove the Tc4331 / Tc43311 test case rows that were added (rows 41-43).
$ws1.Range("A41:J43").Clear()

# Move the "More tests " sheet's content (B3:B7) back onto the main sheet,
# landing at B42:B46 - this also carries over cell formatting (bold+wrap
# style) and the taller row height on the last row.
$ws2.Range("B3:B7").Copy()
$ws1.Range("B42").PasteSpecial(-4104) # xlPasteAll
$excel.CutCopyMode = $false

# The "More tests " sheet is now redundant - delete it entirely.
$ws2.Delete()

# Restore the selection/active cell to where it ends up after the edit.
$ws1.Activate()
$ws1.Range("B45").Select()
